$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N2").Value = "2019-12-31 00:00:00"
$ws.Range("O2").Value = 849350128.1799999
$ws.Range("P2").Value = 42887232439.68
$ws.Range("Q2").Value = 41521780358.8
$ws.Range("R2").Value = 8.1823417591
$ws.Range("S2").Value = 39109365895.91
$ws.Range("T2").Value = 39109365895.91
$ws.Range("U2").Value = 8.3551428659
$ws.Range("V2").Value = 198948372.4
$ws.Range("W2").Value = 1213474541.84
$ws.Range("X2").Value = 684381312.9299999
$ws.Range("Y2").Value = 1361235026.37
$ws.Range("Z2").Value = 1305562660.47
$ws.Range("AA2").Value = 196469033.84
$ws.Range("AG2").Value = 171082795.54
$ws.Range("AP2").Value = 8.2927189835
$ws.Range("AQ2").Value = 1.195763075341
$ws.Range("AR2").Value = -4.04
$ws.Range("AS2").Value = 823538178.86
$ws.Range("AT2").Value = 11.52183359023
